{"js": "// Replace the 25 \"three-digit \u00f7 one-digit\" practice answers in the\n// single table of the document. Each data row of the table (rows\n// 0, 4, 8, 12, 16 \u2014 the other rows are blank \"work space\" rows) holds\n// 5 answer cells; we update every one of them from old text to new\n// text, matched by exact (row, column) position so that duplicate\n// values elsewhere in the table can never cause an ambiguous match.\nconst replacements = [\n  [0, 0, \"551\u00f74=137, 3\", \"534\u00f75=106, 4\"],\n  [0, 1, \"281\u00f77=40, 1\", \"919\u00f77=131, 2\"],\n  [0, 2, \"199\u00f78=24, 7\", \"646\u00f74=161, 2\"],\n  [0, 3, \"147\u00f77=21, 0\", \"350\u00f75=70, 0\"],\n  [0, 4, \"919\u00f72=459, 1\", \"227\u00f72=113, 1\"],\n  [4, 0, \"992\u00f73=330, 2\", \"407\u00f78=50, 7\"],\n  [4, 1, \"270\u00f75=54, 0\", \"273\u00f79=30, 3\"],\n  [4, 2, \"990\u00f78=123, 6\", \"386\u00f76=64, 2\"],\n  [4, 3, \"466\u00f77=66, 4\", \"262\u00f77=37, 3\"],\n  [4, 4, \"590\u00f72=295, 0\", \"245\u00f76=40, 5\"],\n  [8, 0, \"618\u00f75=123, 3\", \"921\u00f78=115, 1\"],\n  [8, 1, \"688\u00f75=137, 3\", \"320\u00f78=40, 0\"],\n  [8, 2, \"598\u00f74=149, 2\", \"664\u00f77=94, 6\"],\n  [8, 3, \"495\u00f77=70, 5\", \"211\u00f74=52, 3\"],\n  [8, 4, \"222\u00f75=44, 2\", \"423\u00f74=105, 3\"],\n  [12, 0, \"122\u00f75=24, 2\", \"928\u00f78=116, 0\"],\n  [12, 1, \"468\u00f73=156, 0\", \"188\u00f78=23, 4\"],\n  [12, 2, \"661\u00f72=330, 1\", \"180\u00f72=90, 0\"],\n  [12, 3, \"627\u00f74=156, 3\", \"497\u00f76=82, 5\"],\n  [12, 4, \"237\u00f77=33, 6\", \"479\u00f78=59, 7\"],\n  [16, 0, \"206\u00f77=29, 3\", \"444\u00f77=63, 3\"],\n  [16, 1, \"732\u00f74=183, 0\", \"230\u00f74=57, 2\"],\n  [16, 2, \"158\u00f74=39, 2\", \"336\u00f76=56, 0\"],\n  [16, 3, \"497\u00f76=82, 5\", \"883\u00f76=147, 1\"],\n  [16, 4, \"537\u00f72=268, 1\", \"679\u00f78=84, 7\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nfor (const [row, col, oldText, newText] of replacements) {\n  const cell = table.getCell(row, col);\n  const results = cell.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, \"Replace\");\n  } else {\n    // Fallback: if the exact text wasn't found (e.g. already edited),\n    // replace the whole cell body text directly.\n    cell.body.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 \"three-digit \u00f7 one-digit\" practice answers in the\n# single table of the document. Each data row of the table (1-based\n# rows 1, 5, 9, 13, 17 -- the other rows are blank \"work space\" rows)\n# holds 5 answer cells; we update every one of them from old text to\n# new text, matched by exact (row, column) position so duplicate\n# values elsewhere in the table can never cause an ambiguous match.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n  @{Row=1;  Col=1; Old=\"551\u00f74=137, 3\"; New=\"534\u00f75=106, 4\"},\n  @{Row=1;  Col=2; Old=\"281\u00f77=40, 1\";  New=\"919\u00f77=131, 2\"},\n  @{Row=1;  Col=3; Old=\"199\u00f78=24, 7\";  New=\"646\u00f74=161, 2\"},\n  @{Row=1;  Col=4; Old=\"147\u00f77=21, 0\";  New=\"350\u00f75=70, 0\"},\n  @{Row=1;  Col=5; Old=\"919\u00f72=459, 1\"; New=\"227\u00f72=113, 1\"},\n  @{Row=5;  Col=1; Old=\"992\u00f73=330, 2\"; New=\"407\u00f78=50, 7\"},\n  @{Row=5;  Col=2; Old=\"270\u00f75=54, 0\";  New=\"273\u00f79=30, 3\"},\n  @{Row=5;  Col=3; Old=\"990\u00f78=123, 6\"; New=\"386\u00f76=64, 2\"},\n  @{Row=5;  Col=4; Old=\"466\u00f77=66, 4\";  New=\"262\u00f77=37, 3\"},\n  @{Row=5;  Col=5; Old=\"590\u00f72=295, 0\"; New=\"245\u00f76=40, 5\"},\n  @{Row=9;  Col=1; Old=\"618\u00f75=123, 3\"; New=\"921\u00f78=115, 1\"},\n  @{Row=9;  Col=2; Old=\"688\u00f75=137, 3\"; New=\"320\u00f78=40, 0\"},\n  @{Row=9;  Col=3; Old=\"598\u00f74=149, 2\"; New=\"664\u00f77=94, 6\"},\n  @{Row=9;  Col=4; Old=\"495\u00f77=70, 5\";  New=\"211\u00f74=52, 3\"},\n  @{Row=9;  Col=5; Old=\"222\u00f75=44, 2\";  New=\"423\u00f74=105, 3\"},\n  @{Row=13; Col=1; Old=\"122\u00f75=24, 2\";  New=\"928\u00f78=116, 0\"},\n  @{Row=13; Col=2; Old=\"468\u00f73=156, 0\"; New=\"188\u00f78=23, 4\"},\n  @{Row=13; Col=3; Old=\"661\u00f72=330, 1\"; New=\"180\u00f72=90, 0\"},\n  @{Row=13; Col=4; Old=\"627\u00f74=156, 3\"; New=\"497\u00f76=82, 5\"},\n  @{Row=13; Col=5; Old=\"237\u00f77=33, 6\";  New=\"479\u00f78=59, 7\"},\n  @{Row=17; Col=1; Old=\"206\u00f77=29, 3\";  New=\"444\u00f77=63, 3\"},\n  @{Row=17; Col=2; Old=\"732\u00f74=183, 0\"; New=\"230\u00f74=57, 2\"},\n  @{Row=17; Col=3; Old=\"158\u00f74=39, 2\";  New=\"336\u00f76=56, 0\"},\n  @{Row=17; Col=4; Old=\"497\u00f76=82, 5\";  New=\"883\u00f76=147, 1\"},\n  @{Row=17; Col=5; Old=\"537\u00f72=268, 1\"; New=\"679\u00f78=84, 7\"}\n)\n\nforeach ($item in $replacements) {\n  $cell = $t.Cell($item.Row, $item.Col)\n  $cell.Range.Text = $item.New\n}\n"}
